# load from file #6 fix
# Replace the process-snapshot data (Name/PID/Memory) with the updated
# snapshot, which has 62 data rows (one more than before) and a fully
# re-sorted Memory column. Also extend the line chart series range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    "chrome.exe",
    "netbeans64.exe",
    "svchost.exe",
    "java.exe",
    "MsMpEng.exe",
    "Skype.exe",
    "SkypeBrowserHost.exe",
    "ekrn.exe",
    "explorer.exe",
    "SearchIndexer.exe",
    "ExpressTray.exe",
    "dwm.exe",
    "Garmin.Cartography.MapUpdate.CoreService.exe",
    "egui.exe",
    "IAStorDataMgrSvc.exe",
    "IAStorIcon.exe",
    "taskhostex.exe",
    "OSPPSVC.EXE",
    "RuntimeBroker.exe",
    "csrss.exe",
    "lsass.exe",
    "audiodg.exe",
    "spd.exe",
    "WmiPrvSE.exe",
    "cfosspeed.exe",
    "dasHost.exe",
    "taskhost.exe",
    "MpCmdRun.exe",
    "4game-service.exe",
    "CCC.exe",
    "iSCTAgent.exe",
    "RAVCpl64.exe",
    "conhost.exe",
    "jusched.exe",
    "atieclxx.exe",
    "spoolsv.exe",
    "MOM.exe",
    "services.exe",
    "LMS.exe",
    "tasklist.exe",
    "mdm.exe",
    "SearchProtocolHost.exe",
    "winlogon.exe",
    "SearchFilterHost.exe",
    "slimsvc.exe",
    "WUDFHost.exe",
    "iSCTsysTray8.exe",
    "sqlwriter.exe",
    "IPROSetMonitor.exe",
    "IOMonitorSrv.exe",
    "jhi_service.exe",
    "SbieSvc.exe",
    "atiesrxx.exe",
    "IntelMeFWService.exe",
    "wininit.exe",
    "Start8_64.exe",
    "NisSrv.exe",
    "System",
    "smss.exe",
    "wmpnetwk.exe",
    "Start8Srv.exe",
    "System Idle Process"
)

# PID column is stored as text in the workbook (e.g. "4904"), not as a
# number -- keep it that way.
$pids = @(
    "4904",
    "4560",
    "756",
    "732",
    "2260",
    "2116",
    "5728",
    "908",
    "3256",
    "3764",
    "4220",
    "884",
    "1748",
    "3296",
    "3960",
    "4072",
    "3304",
    "8100",
    "2196",
    "496",
    "684",
    "7636",
    "1588",
    "4088",
    "584",
    "2968",
    "5768",
    "4800",
    "1484",
    "4684",
    "1220",
    "4112",
    "1040",
    "4328",
    "2808",
    "1308",
    "4320",
    "656",
    "3532",
    "4796",
    "1604",
    "7504",
    "648",
    "3112",
    "1616",
    "2980",
    "876",
    "2168",
    "1232",
    "1568",
    "1472",
    "1032",
    "964",
    "1600",
    "572",
    "1156",
    "4400",
    "4",
    "344",
    "4284",
    "1124",
    "0"
)

$mems = @(
    4263844.0,
    1622600.0,
    274340.0,
    162572.0,
    135400.0,
    129888.0,
    121872.0,
    101536.0,
    81228.0,
    39128.0,
    36272.0,
    32864.0,
    31460.0,
    25380.0,
    24832.0,
    17160.0,
    15424.0,
    11868.0,
    11624.0,
    11612.0,
    11376.0,
    10732.0,
    10368.0,
    10064.0,
    9360.0,
    9092.0,
    9060.0,
    8644.0,
    8340.0,
    7780.0,
    7764.0,
    7060.0,
    7004.0,
    6888.0,
    6292.0,
    6072.0,
    6008.0,
    5956.0,
    5876.0,
    5816.0,
    5756.0,
    5104.0,
    5080.0,
    4988.0,
    4604.0,
    4376.0,
    4224.0,
    4184.0,
    4112.0,
    3580.0,
    3428.0,
    3380.0,
    3340.0,
    3124.0,
    3116.0,
    3020.0,
    2624.0,
    2492.0,
    848.0,
    572.0,
    528.0,
    4.0
)

for ($i = 0; $i -lt $names.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]

    # Force the PID cell to text so a value like "4904" does not get
    # auto-converted to a number; ClearFormats afterwards drops the
    # temporary "@" number format so no stray style sticks around.
    $pidCell = $ws.Cells.Item($row, 2)
    $pidCell.NumberFormat = "@"
    $pidCell.Value = $pids[$i]
    $pidCell.ClearFormats()

    $ws.Cells.Item($row, 3).Value = $mems[$i]
}

# The data now spans rows 2-63 (was 2-62) -- widen the chart series to match.
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$series = $chart.SeriesCollection(1)
$series.Formula = "=SERIES(,Sheet1!`$A`$2:`$A`$63,Sheet1!`$C`$2:`$C`$63,1)"

